$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing header labels
$ws.Range("I1").Value = "Electric Choice ID"
$ws.Range("J1").Value = "Electric Rate Code"
$ws.Range("L1").Value = "Electric Usage (kWh)"

# Add new header cells, copying the style/format of an existing header cell
$ws.Range("K1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("M1").Value = "Gas Choice ID"
$ws.Range("N1").Value = "Gas Rate Code"
$ws.Range("O1").Value = "Gas Usage (therms)"

# Fill in row 2 data
$ws.Range("K2").Value = "N/A"
$ws.Range("M2").Value = "N/A"
$ws.Range("N2").Value = "N/A"
$ws.Range("O2").Value = "N/A"
